$wb = $excel.ActiveWorkbook

# --- Sheet "Victimas" (sheet1) ---
$wsVictimas = $wb.Worksheets.Item("Victimas")
$wsVictimas.Range("A8").Value = 2025
$wsVictimas.Range("B8").Value = 1263
$wsVictimas.Range("C8").Value = 470

# --- Sheet "Ofensores" (sheet2) ---
$wsOfensores = $wb.Worksheets.Item("Ofensores")
$wsOfensores.Range("A8").Value = 2025
$wsOfensores.Range("B8").Value = 320
$wsOfensores.Range("C8").Value = 1381

# Select a cell on each sheet to mirror recorded selections, then make
# "Ofensores" the active sheet/tab as in the final saved state.
$wsVictimas.Range("C12").Select()
$wsOfensores.Range("E5").Select()
$wsOfensores.Activate()
